$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Matlab Postburial" (rId1) - data corrections
# ------------------------------------------------------------------
$wsPostburial = $wb.Worksheets.Item("Matlab Postburial")
$wsPostburial.Range("AH2").Value = 80.290000000000006
$wsPostburial.Range("AI2").Value = 8.0500000000000007
$wsPostburial.Range("AH3").Value = 56.81
$wsPostburial.Range("AI3").Value = 7.88
$wsPostburial.Range("AH4").Value = 61.67

# ------------------------------------------------------------------
# Sheet "Matlab Erate corrected model1" (rId3) - data corrections
# ------------------------------------------------------------------
$wsModel1 = $wb.Worksheets.Item("Matlab Erate corrected model1")

$wsModel1.Range("J2").Value = 29794
$wsModel1.Range("Y2").Value = 5587.0429566990088
$wsModel1.Range("Y2").NumberFormat = "0"

$wsModel1.Range("J3").Value = 25798
$wsModel1.Range("Y3").Value = 6943.8444683042835
$wsModel1.Range("Y3").NumberFormat = "0"

$wsModel1.Range("J4").Value = 10686
$wsModel1.Range("Y4").Value = 3628.9508125627717
$wsModel1.Range("Y4").NumberFormat = "0"

$wsModel1.Range("J5").Value = 31815.169288005323
$wsModel1.Range("Y5").Value = 2633.2125779328203

$wsModel1.Range("J6").Value = 15259.426325440414
$wsModel1.Range("Y6").Value = 2388.2039990666008

$wsModel1.Range("J7").Value = 57967.734988344331
$wsModel1.Range("Y7").Value = 2632.6109411047437

$wsModel1.Range("J8").Value = 24574.288922747248
$wsModel1.Range("Y8").Value = 1859.7831260082412

# ------------------------------------------------------------------
# Sheet "Matlab Erate corrected model2" (rId4) - data corrections
# ------------------------------------------------------------------
$wsModel2 = $wb.Worksheets.Item("Matlab Erate corrected model2")

$wsModel2.Range("J2").Value = 30460
$wsModel2.Range("Y2").Value = 5389.2597821964382
$wsModel2.Range("Y2").NumberFormat = "0"

$wsModel2.Range("J3").Value = 33206.169288005323
$wsModel2.Range("Y3").Value = 2274.8053720227599

$wsModel2.Range("J4").Value = 16551.426325440414
$wsModel2.Range("Y4").Value = 2177.8003446500106

$wsModel2.Range("J5").Value = 58822.734988344331
$wsModel2.Range("Y5").Value = 2540.0402294499991

$wsModel2.Range("J6").Value = 25481.288922747248
$wsModel2.Range("Y6").Value = 1702.0415023685484

# ------------------------------------------------------------------
# View state: selections on each sheet, and which tab is active.
# Set the (soon to be) non-active sheets' selections first, then
# activate + select on the sheet that should end up active so that
# tabSelected / activeTab line up with the final active sheet.
# ------------------------------------------------------------------
$wsPostburial.Activate() | Out-Null
$wsPostburial.Range("B16").Select() | Out-Null

$wsModel1.Activate() | Out-Null
$wsModel1.Range("AA20").Select() | Out-Null

$wsModel2.Activate() | Out-Null
$wsModel2.Range("AB24").Select() | Out-Null
